$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 'producto' as a new header column (P1), matching the look of
#     the other header cells (bold, centered, bordered) by copying O1's
#     formatting onto the new cell.
$ws.Range("P1").Value = "producto"
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Small fix: correct the test client's fecha_ingreso (E2) by one day.
#     The source value is a plain text string ("yyyy-mm-dd"), not a real
#     date, so force text storage before assigning to avoid Excel
#     auto-converting it to a date serial number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2025-12-11"
$ws.Range("E2").Style = "Normal"      # drop the temporary text format again
